$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").Value = "'96.670.09"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.53%  "

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").Value = "'3.678.76"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.49%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Range("E4").Value = "  -0.15%  "

# Row 5: Solana -> Solana
$ws.Range("D5").Value = "'239.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.00%  "

# Row 6: XRP -> XRP
$ws.Range("E6").Value = "  +7.70%  "

# Row 7: BNB -> BNB
$ws.Range("D7").Value = "'656.35"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.55%  "

# Row 8: Dogecoin -> Dogecoin
$ws.Range("E8").Value = "  +0.42%  "

# Row 9: Cardano -> Cardano
$ws.Range("E9").Value = "  +0.23%  "

# Row 10: USDC -> USDC
$ws.Range("D10").Value = "'0.999"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.04%  "

# Row 11: LidoStakedEther -> LidoStakedEther
$ws.Range("D11").Value = "'3.676.74"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.55%  "

# Row 12: Avalanche -> Avalanche
$ws.Range("D12").Value = "'45.48"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.01%  "

# Row 14: Toncoin -> Toncoin
$ws.Range("D14").Value = "'6.83"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.99%  "

# Row 15: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'4.363.12"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.49%  "

# Row 16: ShibaInu -> ShibaInu
$ws.Range("E16").Value = "  +2.94%  "

# Row 17: WrappedBTC -> WrappedBTC
$ws.Range("D17").Value = "'96.493.62"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.54%  "

# Row 18: WrappedEther -> WrappedEther
$ws.Range("D18").Value = "'3.671.91"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.48%  "

# Row 19: Chainlink -> Chainlink
$ws.Range("D19").Value = "'18.89"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.23%  "

# Row 20: Uniswap -> Uniswap
$ws.Range("D20").Value = "'12.77"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.11%  "

# Row 21: Polkadot -> Polkadot
$ws.Range("D21").Value = "'7.76"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.32%  "

# Row 22: Stellar -> Stellar
$ws.Range("D22").Value = "'0.524"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.50%  "

# Row 23: BitcoinCash -> BitcoinCash
$ws.Range("D23").Value = "'530.50"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.03%  "

# Row 24: SuiNetwork -> SuiNetwork
$ws.Range("D24").Value = "'3.46"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.33%  "

# Row 25: NEARProtocol -> NEARProtocol
$ws.Range("D25").Value = "'7.11"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.49%  "

# Row 26: PEPE -> PEPE
$ws.Range("E26").Value = "  -2.35%  "

# Row 27: Litecoin -> Litecoin
$ws.Range("D27").Value = "'101.77"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.87%  "

# Row 28: Aptos -> Aptos
$ws.Range("D28").Value = "'13.21"
$ws.Range("D28").ClearFormats()

# Row 29: WrappedeETH -> Hedera
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.168"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +8.63%  "

# Row 30: Hedera -> InternetComputer(DFINITY)
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'12.49"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.73%  "

# Row 31: InternetComputer(DFINITY) -> PancakeSwap
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'3.04"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.56%  "

# Row 32: PancakeSwap -> Dai
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.14%  "

# Row 33: Dai -> Fetch.AI
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.91"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +16.26%  "

# Row 34: Fetch.AI -> Cronos
$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D34").Value = "'0.186"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.47%  "

# Row 35: Cronos -> Bittensor
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "'679.29"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +10.82%  "

# Row 36: Bittensor -> Binance-PegBSC-USD
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.61%  "

# Row 37: Binance-PegBSC-USD -> EthereumClassic
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "'32.55"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.97%  "

# Row 38: EthereumClassic -> PolygonEcosystemToken
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.592"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.50%  "

# Row 39: PolygonEcosystemToken -> RenderToken
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'8.85"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.65%  "

# Row 40: RenderToken -> Kaspa
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.160"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.30%  "

# Row 41: Kaspa -> ImmutableX
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Value = "'2.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.27%  "

# Row 42: ImmutableX -> Filecoin
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'6.60"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +9.95%  "

# Row 43: ARBITRUM -> EnergySwap
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'39.73"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +20.33%  "

# Row 44: Filecoin -> ARBITRUM
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'0.961"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.32%  "

# Row 45: EnergySwap -> USDe
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.07%  "

# Row 46: USDe -> VeChain
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0466"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.77%  "

# Row 47: VeChain -> Algorand
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.442"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +10.83%  "

# Row 48: Algorand -> Stacks
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'2.32"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.71%  "

# Row 49: Stacks -> MantraDAO
$ws.Range("B49").Value = "MantraDAO"
$ws.Range("C49").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D49").Value = "'3.72"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.42%  "

# Row 50: MantraDAO -> WhiteBITCoin
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "'23.66"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.10%  "

# Row 51: WhiteBITCoin -> Cosmos
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'8.62"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.06%  "

Write-Host "Applied cryptos update"
